$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column A ("type") for existing rows 2-16: every row becomes "critical" ---
# (rows that were already "critical" are simply re-set to the same value)
$criticalRows = 2,3,4,5,6,7,8,9,10,11,12,13,14,15,16
foreach ($r in $criticalRows) {
    $ws.Cells.Item($r, 1).Value = "critical"
}

# --- Append the new block of rows 17-31: same description/duration pairs as
#     rows 2-16 above, but all tagged "regular" ---
$data = @(
    @{ Desc = "There are no lights or indications of power"; Dur = 10 },
    @{ Desc = "Screen freezes"; Dur = 20 },
    @{ Desc = "Junction box that is uncovered"; Dur = 40 },
    @{ Desc = "Flickering light"; Dur = 30 },
    @{ Desc = "Switches of light not working"; Dur = 10 },
    @{ Desc = "Turned off randomly"; Dur = 50 },
    @{ Desc = "Remote doesn’t work properly"; Dur = 20 },
    @{ Desc = "Not heating"; Dur = 20 },
    @{ Desc = "Not cooling"; Dur = 20 },
    @{ Desc = "Weird Smell"; Dur = 20 },
    @{ Desc = "Buttons not working"; Dur = 40 },
    @{ Desc = "Not working"; Dur = 50 },
    @{ Desc = "Leaking"; Dur = 40 },
    @{ Desc = "Turns on by itself"; Dur = 30 },
    @{ Desc = "No lights"; Dur = 30 }
)

$row = 17
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = "regular"
    $ws.Cells.Item($row, 2).Value = $item.Desc
    $ws.Cells.Item($row, 3).Value = $item.Dur
    $row = $row + 1
}

# --- Selection, matching the authored edit ---
$ws.Range("A15:A16").Select()
